$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (2-7) while preserving the header row and its shared strings.
$ws.Range("A2:T7").ClearContents()

# Write the "Sending cluster" (A) and "Target cluster" (D) columns first, in an order
# that introduces each distinct cluster name (ECs, FAPs, MuSCs) exactly once, in the
# sequence the refreshed NATMI export uses, before the repeated Ligand/Receptor symbols.
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("A5").Value = "FAPs"
$ws.Range("A6").Value = "MuSCs"
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "MuSCs"

# Ligand / Receptor symbol columns (constant across all rows in this sheet).
$ws.Range("B2").Value = "Ccl28"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("B3").Value = "Ccl28"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("B4").Value = "Ccl28"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("B5").Value = "Ccl28"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("B6").Value = "Ccl28"
$ws.Range("C6").Value = "Ccr10"
$ws.Range("B7").Value = "Ccl28"
$ws.Range("C7").Value = "Ccr10"

# Recomputed NATMI TPM metrics for each row.
# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.05794866666666667
$ws.Range("H2").Value = 0.173846
$ws.Range("I2").Value = 0.235800444619869
$ws.Range("J2").Value = 0.235800444619869
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.238415
$ws.Range("N2").Value = 3.715245
$ws.Range("O2").Value = 0.4359607654144799
$ws.Range("P2").Value = 0.4359607654144798
$ws.Range("Q2").Value = 0.07176449803000001
$ws.Range("R2").Value = 0.6458804822700001
$ws.Range("S2").Value = 0.1027997423215528
$ws.Range("T2").Value = 0.1027997423215528
# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.05794866666666667
$ws.Range("H3").Value = 0.173846
$ws.Range("I3").Value = 0.235800444619869
$ws.Range("J3").Value = 0.235800444619869
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.602242
$ws.Range("N3").Value = 4.806725999999999
$ws.Range("O3").Value = 0.5640392345855201
$ws.Range("P3").Value = 0.5640392345855201
$ws.Range("Q3").Value = 0.09284778757733332
$ws.Range("R3").Value = 0.8356300881959999
$ws.Range("S3").Value = 0.1330007022983162
$ws.Range("T3").Value = 0.1330007022983162
# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.093703
$ws.Range("H4").Value = 0.281109
$ws.Range("I4").Value = 0.3812893433650861
$ws.Range("J4").Value = 0.3812893433650861
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.238415
$ws.Range("N4").Value = 3.715245
$ws.Range("O4").Value = 0.4359607654144799
$ws.Range("P4").Value = 0.4359607654144798
$ws.Range("Q4").Value = 0.116043200745
$ws.Range("R4").Value = 1.044388806705
$ws.Range("S4").Value = 0.1662271939778274
$ws.Range("T4").Value = 0.1662271939778273
# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.093703
$ws.Range("H5").Value = 0.281109
$ws.Range("I5").Value = 0.3812893433650861
$ws.Range("J5").Value = 0.3812893433650861
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.602242
$ws.Range("N5").Value = 4.806725999999999
$ws.Range("O5").Value = 0.5640392345855201
$ws.Range("P5").Value = 0.5640392345855201
$ws.Range("Q5").Value = 0.150134882126
$ws.Range("R5").Value = 1.351213939134
$ws.Range("S5").Value = 0.2150621493872587
$ws.Range("T5").Value = 0.2150621493872587
# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.09410133333333333
$ws.Range("H6").Value = 0.282304
$ws.Range("I6").Value = 0.3829102120150449
$ws.Range("J6").Value = 0.3829102120150449
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.238415
$ws.Range("N6").Value = 3.715245
$ws.Range("O6").Value = 0.4359607654144799
$ws.Range("P6").Value = 0.4359607654144798
$ws.Range("Q6").Value = 0.11653650272
$ws.Range("R6").Value = 1.04882852448
$ws.Range("S6").Value = 0.1669338291150997
$ws.Range("T6").Value = 0.1669338291150997
# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.09410133333333333
$ws.Range("H7").Value = 0.282304
$ws.Range("I7").Value = 0.3829102120150449
$ws.Range("J7").Value = 0.3829102120150449
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.602242
$ws.Range("N7").Value = 4.806725999999999
$ws.Range("O7").Value = 0.5640392345855201
$ws.Range("P7").Value = 0.5640392345855201
$ws.Range("Q7").Value = 0.1507731085226666
$ws.Range("R7").Value = 1.356957976704
$ws.Range("S7").Value = 0.2159763828999451
$ws.Range("T7").Value = 0.2159763828999452
